$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('Achocolatado em Pó Nescau 550g', 18.79),
    @('Acém em Pedaços Carrefour Aproximadamente 500 g', 0.008),
    @('Alcatra Bovina Carrefour Aproximadamente 400 g', 0.01),
    @('Alface Lisa Carrefour', 1.2),
    @('Arroz Branco Longo-fino Tipo 1 Tio João 2Kg', 15.29),
    @('Açúcar Refinado União 1kg', 4.69),
    @('Banana Nanica Fresca Orgânica  600g', 7.39),
    @('Banana Prata Fischer Turma da Mônica 750g', 5.29),
    @('Batata Monalisa Carrefour Aprox. 600g', 1.99),
    @('Bisteca Suína Congelada Sadia 1 Kg', 23.39),
    @('Camarão Descascado Cozido 35/45 Maris 400 g', 59.99),
    @('Cebola Carrefour Aprox. 500g', 4.99),
    @('Cenoura Único 1Kg', 3.19),
    @('Contra Filé Swift Mais Aprox. 1,5Kg', 64.90000000000001),
    @('Couve Flor Cledson 300 g', 9.99),
    @('Coxão Mole Fracionado à Vácuo Aprox. 1,3 Kg', 48.99),
    @('Feijão Carioca Tipo 1 Kicaldo 1Kg', 5.99),
    @('Filé de Merluza Congelado Planalto 500 g', 23.9),
    @('Filé de Pescada sem Espinha Swift 500 g', 24.59),
    @('Filé de Tilápia Fresco Carrefour 500 g', 62.29),
    @('Lagarto Swift Mais Aprox. 1,5Kg', 43.59),
    @('Limão Siciliano Carrefour KG', 9.890000000000001),
    @('Linguiça Toscana Swift 700 g', 14.9),
    @('Macarrão Instantâneo Nissin Sabor Galinha Caipira 85g', 1.69),
    @('Macarrão de Sêmola com Ovos Espaguete 8 Adria 500g', 4.09),
    @('Mamão Formosa Sabor & Qualidade Aprox. 1,6 Kg', 9.99),
    @('Manga Palmer Carrefour Aprox. 600g', 4.79),
    @('Maçã Gala Carrefour Aprox 600 g', 6.89),
    @('Melancia Premium Carrefour Aprox.  8Kg', 3.49),
    @('Pera Willians Aprox. 500g', 14.99),
    @('Posta Cação Congelado Buona Pesca 500 g', 22.89),
    @('Salsicha Hot Dog Resfriada Aurora Aproximadamente 500 g', 16.49),
    @('Sorvete Napolitano Nestlé 1,5 Litros', 29.99),
    @('Tomate Carmem Carrefour Aprox. 500g', 0.008),
    @('Uva Escura sem Semente Carrefour 500g', 7.99)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
